$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 64 - this shifts the existing rows 64-109
# (and their data) down to 65-110, preserving all of their values/styles.
$ws.Rows(64).Insert()

# Populate the newly inserted row 64 with the new weekly price record.
$ws.Range("A64").Value2 = 3
$ws.Range("B64").Value2 = "Femacal de La Calera"
$ws.Range("C64").Value2 = "Coquimbo"
$ws.Range("D64").Value2 = 44344
$ws.Range("E64").Value2 = 5
$ws.Range("F64").Value2 = "Fruta"
$ws.Range("G64").Value2 = 100101
$ws.Range("H64").Value2 = "Berries"
$ws.Range("I64").Value2 = 100112025
$ws.Range("J64").Value2 = "Frutilla"
$ws.Range("K64").Value2 = "Sin especificar"
$ws.Range("L64").Value2 = "Especial"
$ws.Range("M64").Value2 = 45
$ws.Range("N64").Value2 = 17000
$ws.Range("O64").Value2 = 17000
$ws.Range("P64").Value2 = 17000
$ws.Range("Q64").Value2 = "$/bandeja 7 kilos"
$ws.Range("R64").Value2 = "Provincia de Melipilla"
$ws.Range("S64").Value2 = 2429
$ws.Range("T64").Value2 = 7
